$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '41.781.23'

# Row 3
$ws.Range("D3").Value = '2.468.90'
$ws.Range("E3").Value = '  -0.56%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.57'
$ws.Range("E5").Value = '  +1.50%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.76'
$ws.Range("E6").Value = '  -0.15%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("E9").Value = '  +3.77%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.84'
$ws.Range("E10").Value = '  +0.96%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0848'
$ws.Range("E11").Value = '  +8.74%  '

# Row 12
$ws.Range("E12").Value = '  +0.27%  '

# Row 13
$ws.Range("D13").Value = '2.849.65'
$ws.Range("E13").Value = '  -0.56%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.91'
$ws.Range("E14").Value = '  +1.17%  '

# Row 15
$ws.Range("E15").Value = '  +2.14%  '

# Row 16
$ws.Range("D16").Value = '2.476.48'
$ws.Range("E16").Value = '  -0.26%  '

# Row 17
$ws.Range("E17").Value = '  +3.98%  '

# Row 18
$ws.Range("D18").Value = '41.750.83'
$ws.Range("E18").Value = '  +0.26%  '

# Row 19
$ws.Range("E19").Value = '  +3.16%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0954'
$ws.Range("E20").Value = '  +3.83%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.60'
$ws.Range("E21").Value = '  +4.37%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.14'
$ws.Range("E22").Value = '  +0.90%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.47'
$ws.Range("E23").Value = '  +1.85%  '

# Row 24
$ws.Range("E24").Value = '  +1.12%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.93'
$ws.Range("E25").Value = '  +1.51%  '

# Row 26
$ws.Range("E26").Value = '  +0.06%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.81'
$ws.Range("E27").Value = '  -0.15%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.27'
$ws.Range("E28").Value = '  +1.34%  '

# Row 29
$ws.Range("E29").Value = '  +1.74%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.23'
$ws.Range("E30").Value = '  +0.34%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.24'
$ws.Range("E31").Value = '  +0.80%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.52'
$ws.Range("E32").Value = '  +2.19%  '

# Row 33
$ws.Range("E33").Value = '  +0.15%  '

# Row 34
$ws.Range("E34").Value = '  +1.68%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.50'
$ws.Range("E35").Value = '  +1.81%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.64'
$ws.Range("E36").Value = '  -2.47%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.91'
$ws.Range("E37").Value = '  -1.42%  '

# Row 38
$ws.Range("E38").Value = '  +1.40%  '

# Row 39
$ws.Range("E39").Value = '  +0.57%  '

# Row 40
$ws.Range("E40").Value = '  -1.09%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.01'
$ws.Range("E41").Value = '  -2.33%  '

# Row 42
$ws.Range("E42").Value = '  -0.09%  '

# Row 43
$ws.Range("D43").Value = '1.975.24'
$ws.Range("E43").Value = '  +1.00%  '

# Row 44
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0284'
$ws.Range("E44").Value = '  +0.15%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.93'
$ws.Range("E45").Value = '  -4.46%  '

# Row 46
$ws.Range("E46").Value = '  +0.20%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.05'
$ws.Range("E47").Value = '  +2.17%  '

# Row 48
$ws.Range("D48").Value = '2.702.93'
$ws.Range("E48").Value = '  -0.81%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.01'
$ws.Range("E49").Value = '  +1.07%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.32'
$ws.Range("E50").Value = '  +0.28%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.16'
$ws.Range("E51").Value = '  +0.17%  '

